# Auto-generated script to apply market-data refresh changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 9714.071  # ALC!H74: 9884.462 -> 9714.071
$ws.Cells.Item(74, 9).Value = 8799.700000000001  # ALC!I74: 8944.223 -> 8799.700000000001
$ws.Cells.Item(74, 11).Value = 8799.700000000001  # ALC!K74: 8944.223 -> 8799.700000000001
$ws.Cells.Item(74, 13).Value = -7863.700000000001  # ALC!M74: -8008.223 -> -7863.700000000001
$ws.Cells.Item(77, 8).Value = 9714.071  # ALC!H77: 9884.462 -> 9714.071
$ws.Cells.Item(77, 9).Value = 8799.700000000001  # ALC!I77: 8944.223 -> 8799.700000000001
$ws.Cells.Item(77, 11).Value = 43998.5  # ALC!K77: 44721.115 -> 43998.5
$ws.Cells.Item(77, 13).Value = -39318.5  # ALC!M77: -40041.115 -> -39318.5
$ws.Cells.Item(97, 8).Value = 2168.3845  # ALC!H97: 2182.4167 -> 2168.3845
$ws.Cells.Item(97, 10).Value = 2168.3845  # ALC!J97: 2182.4167 -> 2168.3845
$ws.Cells.Item(97, 12).Value = 6505.1535  # ALC!L97: 6547.250100000001 -> 6505.1535
$ws.Cells.Item(97, 14).Value = -7497.1535  # ALC!N97: -7539.250100000001 -> -7497.1535
$ws.Cells.Item(138, 8).Value = 2399.9194  # ALC!H138: 2374.492 -> 2399.9194
$ws.Cells.Item(138, 9).Value = 1253.1154  # ALC!I138: 1236.2593 -> 1253.1154
$ws.Cells.Item(138, 11).Value = 3759.3462  # ALC!K138: 3708.7779 -> 3759.3462
$ws.Cells.Item(138, 13).Value = 1380.6538  # ALC!M138: 1431.2221 -> 1380.6538

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3869.4778  # ARM!H32: 3870.3667 -> 3869.4778
$ws.Cells.Item(32, 9).Value = 3362.4158  # ARM!I32: 3363.3147 -> 3362.4158
$ws.Cells.Item(32, 11).Value = 3362.4158  # ARM!K32: 3363.3147 -> 3362.4158
$ws.Cells.Item(32, 13).Value = -3075.4158  # ARM!M32: -3076.3147 -> -3075.4158
$ws.Cells.Item(41, 8).Value = 0  # ARM!H41: 15000 -> 0
$ws.Cells.Item(41, 10).Value = 0  # ARM!J41: 15000 -> 0
$ws.Cells.Item(41, 12).Value = 0  # ARM!L41: 15000 -> 0
$ws.Cells.Item(41, 14).Value = $null  # ARM!N41: -15828 -> (removed)
$ws.Cells.Item(70, 8).Value = 0  # ARM!H70: 50577 -> 0
$ws.Cells.Item(70, 10).Value = 0  # ARM!J70: 50577 -> 0
$ws.Cells.Item(70, 12).Value = 0  # ARM!L70: 50577 -> 0
$ws.Cells.Item(70, 14).Value = $null  # ARM!N70: -51117 -> (removed)
$ws.Cells.Item(73, 8).Value = 0  # ARM!H73: 50577 -> 0
$ws.Cells.Item(73, 10).Value = 0  # ARM!J73: 50577 -> 0
$ws.Cells.Item(73, 12).Value = 0  # ARM!L73: 50577 -> 0
$ws.Cells.Item(73, 14).Value = $null  # ARM!N73: -52449 -> (removed)
$ws.Cells.Item(122, 8).Value = 11563.229  # ARM!H122: 11691.914 -> 11563.229
$ws.Cells.Item(122, 9).Value = 2569  # ARM!I122: 2753.6785 -> 2569
$ws.Cells.Item(122, 10).Value = 55035.332  # ARM!J122: 47444.855 -> 55035.332
$ws.Cells.Item(122, 11).Value = 7707  # ARM!K122: 8261.0355 -> 7707
$ws.Cells.Item(122, 12).Value = 165105.996  # ARM!L122: 142334.565 -> 165105.996
$ws.Cells.Item(122, 13).Value = -5257  # ARM!M122: -5811.0355 -> -5257
$ws.Cells.Item(122, 14).Value = -170005.996  # ARM!N122: -147234.565 -> -170005.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 20835546  # BSM!H64: 17545890 -> 20835546
$ws.Cells.Item(64, 9).Value = 33335574  # BSM!I64: 30305102 -> 33335574
$ws.Cells.Item(64, 10).Value = 2166.5  # BSM!J64: 1973.25 -> 2166.5
$ws.Cells.Item(64, 11).Value = 33335574  # BSM!K64: 30305102 -> 33335574
$ws.Cells.Item(64, 12).Value = 2166.5  # BSM!L64: 1973.25 -> 2166.5
$ws.Cells.Item(64, 13).Value = -33335349  # BSM!M64: -30304877 -> -33335349
$ws.Cells.Item(64, 14).Value = -2616.5  # BSM!N64: -2423.25 -> -2616.5
$ws.Cells.Item(67, 8).Value = 20835546  # BSM!H67: 17545890 -> 20835546
$ws.Cells.Item(67, 9).Value = 33335574  # BSM!I67: 30305102 -> 33335574
$ws.Cells.Item(67, 10).Value = 2166.5  # BSM!J67: 1973.25 -> 2166.5
$ws.Cells.Item(67, 11).Value = 33335574  # BSM!K67: 30305102 -> 33335574
$ws.Cells.Item(67, 12).Value = 2166.5  # BSM!L67: 1973.25 -> 2166.5
$ws.Cells.Item(67, 13).Value = -33334794  # BSM!M67: -30304322 -> -33334794
$ws.Cells.Item(67, 14).Value = -3726.5  # BSM!N67: -3533.25 -> -3726.5
$ws.Cells.Item(86, 8).Value = 1885.3334  # BSM!H86: 1975.3636 -> 1885.3334
$ws.Cells.Item(86, 9).Value = 1293.2222  # BSM!I86: 1343 -> 1293.2222
$ws.Cells.Item(86, 11).Value = 1293.2222  # BSM!K86: 1343 -> 1293.2222
$ws.Cells.Item(86, 13).Value = -170.2221999999999  # BSM!M86: -220 -> -170.2221999999999
$ws.Cells.Item(89, 8).Value = 1885.3334  # BSM!H89: 1975.3636 -> 1885.3334
$ws.Cells.Item(89, 9).Value = 1293.2222  # BSM!I89: 1343 -> 1293.2222
$ws.Cells.Item(89, 11).Value = 6466.111  # BSM!K89: 6715 -> 6466.111
$ws.Cells.Item(89, 13).Value = -850.1109999999999  # BSM!M89: -1099 -> -850.1109999999999
$ws.Cells.Item(94, 8).Value = 847.5833  # BSM!H94: 848.0833 -> 847.5833
$ws.Cells.Item(94, 9).Value = 528  # BSM!I94: 528.6667 -> 528
$ws.Cells.Item(94, 11).Value = 528  # BSM!K94: 528.6667 -> 528
$ws.Cells.Item(94, 13).Value = -77  # BSM!M94: -77.66669999999999 -> -77
$ws.Cells.Item(107, 8).Value = 4836.2856  # BSM!H107: 4136.5713 -> 4836.2856
$ws.Cells.Item(107, 9).Value = 3993.2  # BSM!I107: 3992.6667 -> 3993.2
$ws.Cells.Item(107, 10).Value = 6944  # BSM!J107: 5000 -> 6944
$ws.Cells.Item(107, 11).Value = 3993.2  # BSM!K107: 3992.6667 -> 3993.2
$ws.Cells.Item(107, 12).Value = 6944  # BSM!L107: 5000 -> 6944
$ws.Cells.Item(107, 13).Value = -2073.2  # BSM!M107: -2072.6667 -> -2073.2
$ws.Cells.Item(107, 14).Value = -10784  # BSM!N107: -8840 -> -10784
$ws.Cells.Item(134, 8).Value = 8549.6  # BSM!H134: 8787.286 -> 8549.6
$ws.Cells.Item(134, 9).Value = 8498  # BSM!I134: 8749.5 -> 8498
$ws.Cells.Item(134, 11).Value = 25494  # BSM!K134: 26248.5 -> 25494
$ws.Cells.Item(134, 13).Value = -22959  # BSM!M134: -23713.5 -> -22959
$ws.Cells.Item(139, 8).Value = 59999  # BSM!H139: 80449.5 -> 59999
$ws.Cells.Item(139, 10).Value = 0  # BSM!J139: 100900 -> 0
$ws.Cells.Item(139, 12).Value = 0  # BSM!L139: 100900 -> 0
$ws.Cells.Item(139, 14).Value = $null  # BSM!N139: -111180 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1823.4584  # CRP!H16: 1874.3478 -> 1823.4584
$ws.Cells.Item(16, 9).Value = 1533.1765  # CRP!I16: 1535.9412 -> 1533.1765
$ws.Cells.Item(16, 10).Value = 2528.4285  # CRP!J16: 2833.1667 -> 2528.4285
$ws.Cells.Item(16, 11).Value = 1533.1765  # CRP!K16: 1535.9412 -> 1533.1765
$ws.Cells.Item(16, 12).Value = 2528.4285  # CRP!L16: 2833.1667 -> 2528.4285
$ws.Cells.Item(16, 13).Value = -1246.1765  # CRP!M16: -1248.9412 -> -1246.1765
$ws.Cells.Item(16, 14).Value = -3102.4285  # CRP!N16: -3407.1667 -> -3102.4285
$ws.Cells.Item(31, 8).Value = 3041.7646  # CRP!H31: 3214.5862 -> 3041.7646
$ws.Cells.Item(31, 9).Value = 2146.1785  # CRP!I31: 2335.8 -> 2146.1785
$ws.Cells.Item(31, 10).Value = 7221.1665  # CRP!J31: 8707 -> 7221.1665
$ws.Cells.Item(31, 11).Value = 2146.1785  # CRP!K31: 2335.8 -> 2146.1785
$ws.Cells.Item(31, 12).Value = 7221.1665  # CRP!L31: 8707 -> 7221.1665
$ws.Cells.Item(31, 13).Value = -1851.1785  # CRP!M31: -2040.8 -> -1851.1785
$ws.Cells.Item(31, 14).Value = -7811.1665  # CRP!N31: -9297 -> -7811.1665
$ws.Cells.Item(34, 8).Value = 3041.7646  # CRP!H34: 3214.5862 -> 3041.7646
$ws.Cells.Item(34, 9).Value = 2146.1785  # CRP!I34: 2335.8 -> 2146.1785
$ws.Cells.Item(34, 10).Value = 7221.1665  # CRP!J34: 8707 -> 7221.1665
$ws.Cells.Item(34, 11).Value = 2146.1785  # CRP!K34: 2335.8 -> 2146.1785
$ws.Cells.Item(34, 12).Value = 7221.1665  # CRP!L34: 8707 -> 7221.1665
$ws.Cells.Item(34, 13).Value = -1944.1785  # CRP!M34: -2133.8 -> -1944.1785
$ws.Cells.Item(34, 14).Value = -7625.1665  # CRP!N34: -9111 -> -7625.1665
$ws.Cells.Item(93, 8).Value = 13249.125  # CRP!H93: 16924.25 -> 13249.125
$ws.Cells.Item(93, 9).Value = 13899.286  # CRP!I93: 16924.25 -> 13899.286
$ws.Cells.Item(93, 10).Value = 8698  # CRP!J93: 0 -> 8698
$ws.Cells.Item(93, 11).Value = 13899.286  # CRP!K93: 16924.25 -> 13899.286
$ws.Cells.Item(93, 12).Value = 8698  # CRP!L93: 0 -> 8698
$ws.Cells.Item(93, 13).Value = -12027.286  # CRP!M93: -15052.25 -> -12027.286
$ws.Cells.Item(93, 14).Value = -12442  # CRP!N93: None -> -12442
$ws.Cells.Item(103, 8).Value = 29433.375  # CRP!H103: 36920.168 -> 29433.375
$ws.Cells.Item(103, 9).Value = 27209.572  # CRP!I103: 35304.2 -> 27209.572
$ws.Cells.Item(103, 11).Value = 27209.572  # CRP!K103: 35304.2 -> 27209.572
$ws.Cells.Item(103, 13).Value = -26037.572  # CRP!M103: -34132.2 -> -26037.572
$ws.Cells.Item(113, 8).Value = 1823.4584  # CRP!H113: 1874.3478 -> 1823.4584
$ws.Cells.Item(113, 9).Value = 1533.1765  # CRP!I113: 1535.9412 -> 1533.1765
$ws.Cells.Item(113, 10).Value = 2528.4285  # CRP!J113: 2833.1667 -> 2528.4285
$ws.Cells.Item(113, 11).Value = 1533.1765  # CRP!K113: 1535.9412 -> 1533.1765
$ws.Cells.Item(113, 12).Value = 2528.4285  # CRP!L113: 2833.1667 -> 2528.4285
$ws.Cells.Item(113, 13).Value = 636.8235  # CRP!M113: 634.0588 -> 636.8235
$ws.Cells.Item(113, 14).Value = -6868.4285  # CRP!N113: -7173.1667 -> -6868.4285
$ws.Cells.Item(132, 8).Value = 5822.778  # CRP!H132: 5822.9165 -> 5822.778
$ws.Cells.Item(132, 9).Value = 3492.875  # CRP!I132: 3605.652 -> 3492.875
$ws.Cells.Item(132, 10).Value = 10482.583  # CRP!J132: 9745.77 -> 10482.583
$ws.Cells.Item(132, 11).Value = 10478.625  # CRP!K132: 10816.956 -> 10478.625
$ws.Cells.Item(132, 12).Value = 31447.749  # CRP!L132: 29237.31 -> 31447.749
$ws.Cells.Item(132, 13).Value = -7948.625  # CRP!M132: -8286.956 -> -7948.625
$ws.Cells.Item(132, 14).Value = -36507.749  # CRP!N132: -34297.31 -> -36507.749
$ws.Cells.Item(134, 8).Value = 2431.4902  # CRP!H134: 2476.7 -> 2431.4902
$ws.Cells.Item(134, 9).Value = 2249.2327  # CRP!I134: 2365.1 -> 2249.2327
$ws.Cells.Item(134, 10).Value = 3411.125  # CRP!J134: 2923.1 -> 3411.125
$ws.Cells.Item(134, 11).Value = 6747.6981  # CRP!K134: 7095.299999999999 -> 6747.6981
$ws.Cells.Item(134, 12).Value = 10233.375  # CRP!L134: 8769.299999999999 -> 10233.375
$ws.Cells.Item(134, 13).Value = -4212.6981  # CRP!M134: -4560.299999999999 -> -4212.6981
$ws.Cells.Item(134, 14).Value = -15303.375  # CRP!N134: -13839.3 -> -15303.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 382.75  # CUL!H16: 385.5 -> 382.75
$ws.Cells.Item(16, 9).Value = 189.5  # CUL!I16: 190 -> 189.5
$ws.Cells.Item(16, 10).Value = 576  # CUL!J16: 450.66666 -> 576
$ws.Cells.Item(16, 11).Value = 568.5  # CUL!K16: 570 -> 568.5
$ws.Cells.Item(16, 12).Value = 1728  # CUL!L16: 1351.99998 -> 1728
$ws.Cells.Item(16, 13).Value = -395.5  # CUL!M16: -397 -> -395.5
$ws.Cells.Item(16, 14).Value = -2074  # CUL!N16: -1697.99998 -> -2074
$ws.Cells.Item(33, 8).Value = 557.6  # CUL!H33: 587.55554 -> 557.6
$ws.Cells.Item(33, 9).Value = 345.66666  # CUL!I33: 357.2 -> 345.66666
$ws.Cells.Item(33, 11).Value = 2073.99996  # CUL!K33: 2143.2 -> 2073.99996
$ws.Cells.Item(33, 13).Value = -1790.99996  # CUL!M33: -1860.2 -> -1790.99996
$ws.Cells.Item(122, 8).Value = 1947.7693  # CUL!H122: 1947.6923 -> 1947.7693
$ws.Cells.Item(122, 9).Value = 1395  # CUL!I122: 1305 -> 1395
$ws.Cells.Item(122, 10).Value = 1993.8334  # CUL!J122: 2064.5454 -> 1993.8334
$ws.Cells.Item(122, 11).Value = 12555  # CUL!K122: 11745 -> 12555
$ws.Cells.Item(122, 12).Value = 17944.5006  # CUL!L122: 18580.9086 -> 17944.5006
$ws.Cells.Item(122, 13).Value = -10105  # CUL!M122: -9295 -> -10105
$ws.Cells.Item(122, 14).Value = -22844.5006  # CUL!N122: -23480.9086 -> -22844.5006
$ws.Cells.Item(132, 8).Value = 2528.3044  # CUL!H132: 2528.7827 -> 2528.3044
$ws.Cells.Item(132, 10).Value = 3097.2666  # CUL!J132: 3098 -> 3097.2666
$ws.Cells.Item(132, 12).Value = 27875.3994  # CUL!L132: 27882 -> 27875.3994
$ws.Cells.Item(132, 14).Value = -32935.39939999999  # CUL!N132: -32942 -> -32935.39939999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 75998  # GSM!H46: 64147.668 -> 75998
$ws.Cells.Item(46, 9).Value = 5000  # GSM!I46: 4443 -> 5000
$ws.Cells.Item(46, 10).Value = 93747.5  # GSM!J46: 94000 -> 93747.5
$ws.Cells.Item(46, 11).Value = 5000  # GSM!K46: 4443 -> 5000
$ws.Cells.Item(46, 12).Value = 93747.5  # GSM!L46: 94000 -> 93747.5
$ws.Cells.Item(46, 13).Value = -4844  # GSM!M46: -4287 -> -4844
$ws.Cells.Item(46, 14).Value = -94059.5  # GSM!N46: -94312 -> -94059.5
$ws.Cells.Item(80, 8).Value = 19206294  # GSM!H80: 10959674 -> 19206294
$ws.Cells.Item(80, 9).Value = 25270244  # GSM!I80: 14383228 -> 25270244
$ws.Cells.Item(80, 10).Value = 3781.8333  # GSM!J80: 4299.2 -> 3781.8333
$ws.Cells.Item(80, 11).Value = 25270244  # GSM!K80: 14383228 -> 25270244
$ws.Cells.Item(80, 12).Value = 3781.8333  # GSM!L80: 4299.2 -> 3781.8333
$ws.Cells.Item(80, 13).Value = -25269246  # GSM!M80: -14382230 -> -25269246
$ws.Cells.Item(80, 14).Value = -5777.8333  # GSM!N80: -6295.2 -> -5777.8333
$ws.Cells.Item(83, 8).Value = 19206294  # GSM!H83: 10959674 -> 19206294
$ws.Cells.Item(83, 9).Value = 25270244  # GSM!I83: 14383228 -> 25270244
$ws.Cells.Item(83, 10).Value = 3781.8333  # GSM!J83: 4299.2 -> 3781.8333
$ws.Cells.Item(83, 11).Value = 126351220  # GSM!K83: 71916140 -> 126351220
$ws.Cells.Item(83, 12).Value = 18909.1665  # GSM!L83: 21496 -> 18909.1665
$ws.Cells.Item(83, 13).Value = -126346228  # GSM!M83: -71911148 -> -126346228
$ws.Cells.Item(83, 14).Value = -28893.1665  # GSM!N83: -31480 -> -28893.1665
$ws.Cells.Item(102, 8).Value = 6023.5386  # GSM!H102: 5813.037 -> 6023.5386
$ws.Cells.Item(102, 9).Value = 3433.6  # GSM!I102: 3286.2856 -> 3433.6
$ws.Cells.Item(102, 11).Value = 3433.6  # GSM!K102: 3286.2856 -> 3433.6
$ws.Cells.Item(102, 13).Value = -1811.6  # GSM!M102: -1664.2856 -> -1811.6
$ws.Cells.Item(110, 8).Value = 0  # GSM!H110: 20000 -> 0
$ws.Cells.Item(110, 9).Value = 0  # GSM!I110: 20000 -> 0
$ws.Cells.Item(110, 11).Value = 0  # GSM!K110: 20000 -> 0
$ws.Cells.Item(110, 13).Value = $null  # GSM!M110: -15910 -> (removed)
$ws.Cells.Item(132, 8).Value = 1306.8572  # GSM!H132: 1314.5 -> 1306.8572
$ws.Cells.Item(132, 9).Value = 1306.8572  # GSM!I132: 1377.6 -> 1306.8572
$ws.Cells.Item(132, 10).Value = 0  # GSM!J132: 999 -> 0
$ws.Cells.Item(132, 11).Value = 3920.5716  # GSM!K132: 4132.799999999999 -> 3920.5716
$ws.Cells.Item(132, 12).Value = 0  # GSM!L132: 2997 -> 0
$ws.Cells.Item(132, 13).Value = -1390.5716  # GSM!M132: -1602.799999999999 -> -1390.5716
$ws.Cells.Item(132, 14).Value = $null  # GSM!N132: -8057 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 15626316  # LTW!H16: 1441.375 -> 15626316
$ws.Cells.Item(16, 9).Value = 20834456  # LTW!I16: 1290.1666 -> 20834456
$ws.Cells.Item(16, 11).Value = 20834456  # LTW!K16: 1290.1666 -> 20834456
$ws.Cells.Item(16, 13).Value = -20834286  # LTW!M16: -1120.1666 -> -20834286
$ws.Cells.Item(46, 8).Value = 2862  # LTW!H46: 3320.8333 -> 2862
$ws.Cells.Item(46, 9).Value = 1065  # LTW!I46: 1037 -> 1065
$ws.Cells.Item(46, 10).Value = 4017.2144  # LTW!J46: 4199.231 -> 4017.2144
$ws.Cells.Item(46, 11).Value = 1065  # LTW!K46: 1037 -> 1065
$ws.Cells.Item(46, 12).Value = 4017.2144  # LTW!L46: 4199.231 -> 4017.2144
$ws.Cells.Item(46, 13).Value = -877  # LTW!M46: -849 -> -877
$ws.Cells.Item(46, 14).Value = -4393.2144  # LTW!N46: -4575.231 -> -4393.2144
$ws.Cells.Item(135, 8).Value = 80333  # LTW!H135: 80374 -> 80333
$ws.Cells.Item(135, 10).Value = 80333  # LTW!J135: 80374 -> 80333
$ws.Cells.Item(135, 12).Value = 80333  # LTW!L135: 80374 -> 80333
$ws.Cells.Item(135, 14).Value = -90473  # LTW!N135: -90514 -> -90473
$ws.Cells.Item(137, 8).Value = 60658.668  # LTW!H137: 61033 -> 60658.668
$ws.Cells.Item(137, 10).Value = 55988  # LTW!J137: 59239.6 -> 55988
$ws.Cells.Item(137, 12).Value = 55988  # LTW!L137: 59239.6 -> 55988
$ws.Cells.Item(137, 14).Value = -66188  # LTW!N137: -69439.60000000001 -> -66188

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1069.174  # WVR!H107: 933.55554 -> 1069.174
$ws.Cells.Item(107, 9).Value = 1049.5  # WVR!I107: 936.75 -> 1049.5
$ws.Cells.Item(107, 10).Value = 1099.7778  # WVR!J107: 928.9091 -> 1099.7778
$ws.Cells.Item(107, 11).Value = 3148.5  # WVR!K107: 2810.25 -> 3148.5
$ws.Cells.Item(107, 12).Value = 3299.3334  # WVR!L107: 2786.7273 -> 3299.3334
$ws.Cells.Item(107, 13).Value = -1228.5  # WVR!M107: -890.25 -> -1228.5
$ws.Cells.Item(107, 14).Value = -7139.3334  # WVR!N107: -6626.7273 -> -7139.3334
$ws.Cells.Item(132, 8).Value = 2636.8462  # WVR!H132: 2957.3635 -> 2636.8462
$ws.Cells.Item(132, 9).Value = 2808.3635  # WVR!I132: 3060.4443 -> 2808.3635
$ws.Cells.Item(132, 10).Value = 1693.5  # WVR!J132: 2493.5 -> 1693.5
$ws.Cells.Item(132, 11).Value = 8425.0905  # WVR!K132: 9181.332900000001 -> 8425.0905
$ws.Cells.Item(132, 12).Value = 5080.5  # WVR!L132: 7480.5 -> 5080.5
$ws.Cells.Item(132, 13).Value = -5895.0905  # WVR!M132: -6651.332900000001 -> -5895.0905
$ws.Cells.Item(132, 14).Value = -10140.5  # WVR!N132: -12540.5 -> -10140.5

Write-Host "Applied all market-data refresh cell updates."